# Fix metadata block parser when Excel first cell is blank (None)
#
# The original sheet starts directly with a table block ("**places_to_go")
# in row 1. This reproduces a bug where the parser chokes if the very
# first cell in the sheet is blank/None, by adding a small metadata block
# (author:/purpose:) and a one-column table ("***read_this_summer") ahead
# of the existing tables, with a blank separator row before each block -
# pushing the existing content down by 7 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert 7 new rows at the top, shifting everything else down.
$ws.Rows("1:7").Insert()

# New metadata block (rows 1-2, row 3 left blank as a separator).
$ws.Range("A1").Value = "author:"
$ws.Range("B1").Value = "XYODA"
$ws.Range("A2").Value = "purpose:"
$ws.Range("B2").Value = "Save the galaxy"

# New one-column table block (rows 4-6, row 7 left blank as a separator).
# (values are interned into the shared-string table in the order they are
# first written, so "War and Peace" is written before "***read_this_summer"
# to reproduce the author's original string order)
$ws.Range("A5").Value = "War and Peace"
$ws.Range("A4").Value = "***read_this_summer"
$ws.Range("A6").Value = "Crime and Punishment"

# Match the author's final selection.
$ws.Range("A7").Select() | Out-Null
